# Auto-generated edit script
# Applies the scheduled-runner price/profit updates to the Leve profit sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 144.77
$ws.Range("I15").Value = 144.77
$ws.Range("K15").Value = 434.3100000000001
$ws.Range("M15").Value = -265.3100000000001
# Row 32
$ws.Range("H32").Value = 1489.1818
$ws.Range("I32").Value = 1200
$ws.Range("J32").Value = 1518.1
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 1518.1
$ws.Range("M32").Value = -874
$ws.Range("N32").Value = -2170.1
# Row 55
$ws.Range("H55").Value = 197.18182
$ws.Range("I55").Value = 99.75
$ws.Range("J55").Value = 252.85715
$ws.Range("K55").Value = 99.75
$ws.Range("L55").Value = 252.85715
$ws.Range("M55").Value = 114.25
$ws.Range("N55").Value = -680.85715
# Row 113
$ws.Range("H113").Value = 2497.6843
$ws.Range("I113").Value = 1995.6957
$ws.Range("J113").Value = 3267.4
$ws.Range("K113").Value = 1995.6957
$ws.Range("L113").Value = 3267.4
$ws.Range("M113").Value = 1258.3043
$ws.Range("N113").Value = -9775.4
# Row 132
$ws.Range("H132").Value = 1450.8
$ws.Range("I132").Value = 1382.0294
$ws.Range("K132").Value = 4146.0882
$ws.Range("M132").Value = -1616.0882
# Row 137
$ws.Range("H137").Value = 1179.9
$ws.Range("I137").Value = 1009.1818
$ws.Range("J137").Value = 1388.5555
$ws.Range("K137").Value = 3027.5454
$ws.Range("L137").Value = 4165.666499999999
$ws.Range("M137").Value = -477.5454
$ws.Range("N137").Value = -9265.666499999999
# Row 140
$ws.Range("H140").Value = 83214.44500000001
$ws.Range("J140").Value = 83214.44500000001
$ws.Range("L140").Value = 83214.44500000001
$ws.Range("N140").Value = -93574.44500000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2634.62
$ws.Range("I32").Value = 2269.7583
$ws.Range("J32").Value = 6323.778
$ws.Range("K32").Value = 2269.7583
$ws.Range("L32").Value = 6323.778
$ws.Range("M32").Value = -1982.7583
$ws.Range("N32").Value = -6897.778
# Row 74
$ws.Range("H74").Value = 3750.5
$ws.Range("I74").Value = 4586.857
$ws.Range("J74").Value = 2914.1428
$ws.Range("K74").Value = 4586.857
$ws.Range("L74").Value = 2914.1428
$ws.Range("M74").Value = -3712.857
$ws.Range("N74").Value = -4662.1428
# Row 77
$ws.Range("H77").Value = 3750.5
$ws.Range("I77").Value = 4586.857
$ws.Range("J77").Value = 2914.1428
$ws.Range("K77").Value = 22934.285
$ws.Range("L77").Value = 14570.714
$ws.Range("M77").Value = -18566.285
$ws.Range("N77").Value = -23306.714
# Row 102
$ws.Range("H102").Value = 6027.3335
$ws.Range("I102").Value = 6338.1
$ws.Range("J102").Value = 5638.875
$ws.Range("K102").Value = 6338.1
$ws.Range("L102").Value = 5638.875
$ws.Range("M102").Value = -4716.1
$ws.Range("N102").Value = -8882.875
# Row 132
$ws.Range("H132").Value = 1927.0193
$ws.Range("I132").Value = 1691.5834
$ws.Range("J132").Value = 2456.75
$ws.Range("K132").Value = 5074.7502
$ws.Range("L132").Value = 7370.25
$ws.Range("M132").Value = -2544.7502
$ws.Range("N132").Value = -12430.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 97
$ws.Range("H97").Value = 5915.5
$ws.Range("I97").Value = 5689.143
$ws.Range("J97").Value = 7500
$ws.Range("K97").Value = 5689.143
$ws.Range("L97").Value = 7500
$ws.Range("M97").Value = -4698.143
$ws.Range("N97").Value = -9482

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 44
$ws.Range("H44").Value = 21200
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 21200
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 21200
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -22084
# Row 86
$ws.Range("H86").Value = 37009.43
$ws.Range("I86").Value = 3586.6667
$ws.Range("J86").Value = 120566.336
$ws.Range("K86").Value = 3586.6667
$ws.Range("L86").Value = 120566.336
$ws.Range("M86").Value = -2463.6667
$ws.Range("N86").Value = -122812.336
# Row 89
$ws.Range("H89").Value = 37009.43
$ws.Range("I89").Value = 3586.6667
$ws.Range("J89").Value = 120566.336
$ws.Range("K89").Value = 17933.3335
$ws.Range("L89").Value = 602831.6799999999
$ws.Range("M89").Value = -12317.3335
$ws.Range("N89").Value = -614063.6799999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Range("H59").Value = 200001220
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 200001220
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 600003660
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -600004740
# Row 92
$ws.Range("H92").Value = 554.6667
$ws.Range("I92").Value = 575.5
$ws.Range("J92").Value = 538
$ws.Range("K92").Value = 1726.5
$ws.Range("L92").Value = 1614
$ws.Range("M92").Value = -478.5
$ws.Range("N92").Value = -4110

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 30160
$ws.Range("I46").Value = 27750
$ws.Range("J46").Value = 39800
$ws.Range("K46").Value = 27750
$ws.Range("L46").Value = 39800
$ws.Range("M46").Value = -27594
$ws.Range("N46").Value = -40112
# Row 49
$ws.Range("H49").Value = 5000
$ws.Range("J49").Value = 5000
$ws.Range("L49").Value = 5000
$ws.Range("N49").Value = -5368
# Row 132
$ws.Range("H132").Value = 2139.4167
$ws.Range("I132").Value = 1758.5
$ws.Range("J132").Value = 3129.8
$ws.Range("K132").Value = 5275.5
$ws.Range("L132").Value = 9389.400000000001
$ws.Range("M132").Value = -2745.5
$ws.Range("N132").Value = -14449.4
# Row 140
$ws.Range("H140").Value = 95300
$ws.Range("J140").Value = 95300
$ws.Range("L140").Value = 95300
$ws.Range("N140").Value = -105660

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 42
$ws.Range("H42").Value = 4900
# Row 46
$ws.Range("H46").Value = 7188.4116
$ws.Range("I46").Value = 1345.5454
$ws.Range("J46").Value = 17900.334
$ws.Range("K46").Value = 1345.5454
$ws.Range("L46").Value = 17900.334
$ws.Range("M46").Value = -1157.5454
$ws.Range("N46").Value = -18276.334
# Row 49
$ws.Range("H49").Value = 4900
# Row 132
$ws.Range("H132").Value = 5119.5654
$ws.Range("I132").Value = 6252.615
$ws.Range("J132").Value = 3646.6
$ws.Range("K132").Value = 18757.845
$ws.Range("L132").Value = 10939.8
$ws.Range("M132").Value = -16227.845
$ws.Range("N132").Value = -15999.8
# Row 136
$ws.Range("H136").Value = 3949.4138
$ws.Range("I136").Value = 2598.3157
$ws.Range("J136").Value = 6516.5
$ws.Range("K136").Value = 7794.9471
$ws.Range("L136").Value = 19549.5
$ws.Range("M136").Value = -5244.9471
$ws.Range("N136").Value = -24649.5
# Row 139
$ws.Range("H139").Value = 54916.668
$ws.Range("J139").Value = 54916.668
$ws.Range("L139").Value = 54916.668
$ws.Range("N139").Value = -65196.668

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 44
$ws.Range("H44").Value = 6000
$ws.Range("J44").Value = 6000
$ws.Range("L44").Value = 6000
$ws.Range("N44").Value = -7108
# Row 141
$ws.Range("H141").Value = 54411.5
$ws.Range("J141").Value = 55123.89
$ws.Range("L141").Value = 55123.89
$ws.Range("M141").Value = -65483.89
$ws.Range("N141").Value = -65483.89
